$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Week 30" column (column AE) header and data.
$ws.Range("AE1").Value = "Week 30"
$ws.Range("AE2").Value = 3.5
$ws.Range("AE3").Value = 2.5
$ws.Range("AE4").Value = 7
$ws.Range("AE5").Value = 10
$ws.Range("AE9").Value = 1.5
